$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "41.723.01"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -0.01%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.475.60"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  +0.06%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "320.76"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +1.51%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "92.12"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.79%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.550"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -0.66%  "
$ws.Range("B10").Value = "Avalanche"
$ws.Range("C10").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "32.99"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0857"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +2.52%  "
$ws.Range("E12").Value = "  -0.93%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "2.857.52"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("E14").Value = "  +0.21%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "15.48"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -2.08%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "2.471.70"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +1.18%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.794"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +1.62%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "41.666.88"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("E19").Value = "  -0.79%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0941"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.45%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "71.27"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("E22").Value = "  -1.47%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "239.52"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.11%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.75"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +1.40%  "
$ws.Range("E25").Value = "  +1.13%  "
$ws.Range("E26").Value = "  +0.07%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "24.99"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +1.53%  "
$ws.Range("E28").Value = "  -1.16%  "
$ws.Range("E29").Value = "  -0.19%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "36.67"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +3.36%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "157.16"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +1.00%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "5.44"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -1.86%  "
$ws.Range("E33").Value = "  +0.02%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.0769"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +1.02%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "2.56"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -0.29%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "17.16"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -1.91%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.84"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +2.42%  "
$ws.Range("E38").Value = "  +1.30%  "
$ws.Range("E39").Value = "  -0.83%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.103"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +0.37%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "4.00"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("E42").Value = "  -2.28%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "2.000.06"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +1.13%  "
$ws.Range("E44").Value = "  +0.28%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "18.64"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -2.08%  "
$ws.Range("E46").Value = "  +0.68%  "
$ws.Range("E47").Value = "  +4.50%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "2.738.45"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +1.32%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "97.73"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.77%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "76.17"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +5.21%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "67.29"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +0.18%  "
